# Replace development_stage_ontology_term_id (column AL) values in the
# "Tier 1_obs" worksheet with the donor's raw age, for rows where the
# ontology term id had previously been filled in as a placeholder.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tier 1_obs")

$ages = @{
    6  = "63"
    7  = "63"
    8  = "63"
    9  = "63"
    10 = "63"
    11 = "64"
    12 = "64"
    13 = "64"
    14 = "64"
    15 = "64"
    16 = "50"
    17 = "50"
    18 = "50"
    19 = "50"
    20 = "50"
    21 = "50"
    22 = "54"
    23 = "54"
    24 = "54"
    25 = "54"
    26 = "54"
    27 = "69"
    28 = "69"
    29 = "69"
    30 = "69"
    31 = "57"
    32 = "57"
    33 = "57"
    34 = "57"
}

foreach ($row in $ages.Keys) {
    $cell = $ws.Range("AL$row")
    # Force the value to be stored as text (donor age in years) rather
    # than being auto-coerced to a number, then drop the temporary
    # "Text" number format so the cell's style/formatting stays as it was.
    $cell.NumberFormat = "@"
    $cell.Value = $ages[$row]
    $cell.Style = "Normal"
}
